$d = $word.ActiveDocument

# Simple header field replacements
$d.Content.Find.Execute("Date: 2024-12-13", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Date: 2024-12-25", 2)

$d.Content.Find.Execute("Time: 20:05", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Time: 18:07", 2)

$d.Content.Find.Execute("Initiated by: ramsha", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Initiated by: Ramsha Khan", 2)

$d.Content.Find.Execute("Minutes Verified by: ramsha", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Minutes Verified by: Ramsha Khan", 2)

# Agenda paragraph
$d.Content.Find.Execute("1) The agenda for the meeting was to discuss the QFO marketing strategies based on the last quarter's performance.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "The agenda for the meeting includes discussing the project progress from last week, reviewing upcoming deadlines, reviewing bugs reported by the QA team, discussing the deployment process, and allowing time for questions and concerns.", 2)

# Resolution paragraph - find the whole paragraph (the run contains a <w:br/> between two <w:t>
# which Word represents as a paragraph-internal line break within the same Range text as a
# vertical-tab character \x0B). Replace the entire paragraph range text.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "1) The resolutions reached at the meeting include allocating more budget*") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "The decisions or resolutions reached at the meeting include allocating more resources to resolve the roadblocks in API integration, prioritizing certain bugs during the next sprint, finalizing the documentation before the official release, and discussing potential improvements for the deployment process."
        break
    }
}

# Summary paragraph
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "1) The discussion focused on the ineffectiveness*") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "The meeting discussed the progress of the project, the roadblocks in API integration, the tight project timelines, the need to prioritize certain bugs, the upcoming release, and the potential improvements for the deployment process."
        break
    }
}
